# Auto-generated script applying scheduled market-data refresh to Leve profit tables
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1622.25
$ws.Range("I86").Value = 1568.2858
$ws.Range("K86").Value = 1568.2858
$ws.Range("M86").Value = -445.2858000000001
$ws.Range("H89").Value = 1622.25
$ws.Range("I89").Value = 1568.2858
$ws.Range("K89").Value = 7841.429
$ws.Range("M89").Value = -2225.429
$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -80
$ws.Range("N107").ClearContents()
$ws.Range("H136").Value = 25028
$ws.Range("J136").Value = 25028
$ws.Range("L136").Value = 25028
$ws.Range("N136").Value = -35228
$ws.Range("H139").Value = 20500.445
$ws.Range("J139").Value = 20500.445
$ws.Range("L139").Value = 20500.445
$ws.Range("N139").Value = -30780.445
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 75532.25
$ws.Range("I53").Value = 2000
$ws.Range("J53").Value = 100043
$ws.Range("K53").Value = 2000
$ws.Range("L53").Value = 100043
$ws.Range("M53").Value = -1318
$ws.Range("N53").Value = -101407
$ws.Range("H74").Value = 2197.9412
$ws.Range("I74").Value = 1740.3572
$ws.Range("J74").Value = 4333.3335
$ws.Range("K74").Value = 1740.3572
$ws.Range("L74").Value = 4333.3335
$ws.Range("M74").Value = -866.3571999999999
$ws.Range("N74").Value = -6081.3335
$ws.Range("H77").Value = 2197.9412
$ws.Range("I77").Value = 1740.3572
$ws.Range("J77").Value = 4333.3335
$ws.Range("K77").Value = 8701.786
$ws.Range("L77").Value = 21666.6675
$ws.Range("M77").Value = -4333.786
$ws.Range("N77").Value = -30402.6675
$ws.Range("H139").Value = 18915.072
$ws.Range("J139").Value = 18915.072
$ws.Range("L139").Value = 18915.072
$ws.Range("N139").Value = -29195.072
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 87527
$ws.Range("J38").Value = 87527
$ws.Range("L38").Value = 87527
$ws.Range("N38").Value = -88359
$ws.Range("H107").Value = 3475
$ws.Range("I107").Value = 1962.5
$ws.Range("J107").Value = 6500
$ws.Range("K107").Value = 1962.5
$ws.Range("L107").Value = 6500
$ws.Range("M107").Value = -42.5
$ws.Range("N107").Value = -10340
$ws.Range("H135").Value = 18378.703
$ws.Range("J135").Value = 18178.303
$ws.Range("L135").Value = 18178.303
$ws.Range("N135").Value = -28318.303
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4911.4287
$ws.Range("I94").Value = 9606
$ws.Range("J94").Value = 3033.6
$ws.Range("K94").Value = 9606
$ws.Range("L94").Value = 3033.6
$ws.Range("M94").Value = -9155
$ws.Range("N94").Value = -3935.6
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1462.2
$ws.Range("I75").Value = 802.5
$ws.Range("J75").Value = 1902
$ws.Range("K75").Value = 2407.5
$ws.Range("L75").Value = 5706
$ws.Range("M75").Value = -1409.5
$ws.Range("N75").Value = -7702
$ws.Range("H78").Value = 1462.2
$ws.Range("I78").Value = 802.5
$ws.Range("J78").Value = 1902
$ws.Range("K78").Value = 7222.5
$ws.Range("L78").Value = 17118
$ws.Range("M78").Value = -2230.5
$ws.Range("N78").Value = -27102
$ws.Range("H94").Value = 3535.125
$ws.Range("I94").Value = 1856
$ws.Range("J94").Value = 3775
$ws.Range("K94").Value = 5568
$ws.Range("L94").Value = 11325
$ws.Range("M94").Value = -4892
$ws.Range("N94").Value = -12677
$ws.Range("H95").Value = 4163.3335
$ws.Range("J95").Value = 3996
$ws.Range("L95").Value = 11988
$ws.Range("N95").Value = -16106
$ws.Range("H96").Value = 3100
$ws.Range("I96").Value = 1750
$ws.Range("J96").Value = 4000
$ws.Range("K96").Value = 5250
$ws.Range("L96").Value = 12000
$ws.Range("M96").Value = -3191
$ws.Range("N96").Value = -16118
$ws.Range("H113").Value = 1283.7693
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1283.7693
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3851.3079
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8191.3079
$ws.Range("H131").Value = 1693.4865
$ws.Range("J131").Value = 1221.3334
$ws.Range("L131").Value = 3664.0002
$ws.Range("N131").Value = -13744.0002
$ws.Range("H137").Value = 3216
$ws.Range("I137").Value = 2522.6667
$ws.Range("J137").Value = 5296
$ws.Range("K137").Value = 7568.000100000001
$ws.Range("L137").Value = 15888
$ws.Range("M137").Value = -2468.000100000001
$ws.Range("N137").Value = -26088
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4366772.5
$ws.Range("I11").Value = 6403400
$ws.Range("J11").Value = 2669583.2
$ws.Range("K11").Value = 6403400
$ws.Range("L11").Value = 2669583.2
$ws.Range("M11").Value = -6403261
$ws.Range("N11").Value = -2669861.2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H15").Value = 29333.334
$ws.Range("J15").Value = 29333.334
$ws.Range("L15").Value = 29333.334
$ws.Range("N15").Value = -29909.334
$ws.Range("H20").Value = 54005.5
$ws.Range("I20").Value = 8005
$ws.Range("J20").Value = 100006
$ws.Range("K20").Value = 8005
$ws.Range("L20").Value = 100006
$ws.Range("M20").Value = -7760
$ws.Range("N20").Value = -100496
$ws.Range("H21").Value = 30007
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 30007
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 30007
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -30353
$ws.Range("H24").Value = 37206
$ws.Range("I24").Value = 69335.336
$ws.Range("J24").Value = 29173.666
$ws.Range("K24").Value = 69335.336
$ws.Range("L24").Value = 29173.666
$ws.Range("M24").Value = -69162.336
$ws.Range("N24").Value = -29519.666
$ws.Range("H27").Value = 15000
$ws.Range("J27").Value = 15000
$ws.Range("L27").Value = 15000
$ws.Range("N27").Value = -15332
$ws.Range("H30").Value = 30007
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 30007
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 30007
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -30217
$ws.Range("H33").Value = 82156.42999999999
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 82156.42999999999
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 82156.42999999999
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = -82660.42999999999
$ws.Range("H81").Value = 29333.334
$ws.Range("J81").Value = 29333.334
$ws.Range("L81").Value = 29333.334
$ws.Range("N81").Value = -31329.334
$ws.Range("H84").Value = 29333.334
$ws.Range("J84").Value = 29333.334
$ws.Range("L84").Value = 88000.00199999999
$ws.Range("N84").Value = -97984.00199999999
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 142862140
$ws.Range("I61").Value = 250000750
$ws.Range("J61").Value = 10666.667
$ws.Range("K61").Value = 250000750
$ws.Range("L61").Value = 10666.667
$ws.Range("M61").Value = -250000548
$ws.Range("N61").Value = -11070.667
$ws.Range("H113").Value = 142862140
$ws.Range("I113").Value = 250000750
$ws.Range("J113").Value = 10666.667
$ws.Range("K113").Value = 250000750
$ws.Range("L113").Value = 10666.667
$ws.Range("M113").Value = -249998580
$ws.Range("N113").Value = -15006.667
$ws.Range("H135").Value = 46300
$ws.Range("J135").Value = 46300
$ws.Range("L135").Value = 46300
$ws.Range("N135").Value = -56440
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 100007
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 100007
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100007
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -100287
$ws.Range("H14").Value = 1999.5
$ws.Range("I14").Value = 1999.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1999.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1831.5
$ws.Range("N14").ClearContents()
$ws.Range("H19").Value = 20000
$ws.Range("J19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("N19").Value = -20348
$ws.Range("H107").Value = 1535.75
$ws.Range("I107").Value = 326.57144
$ws.Range("K107").Value = 979.71432
$ws.Range("M107").Value = 940.28568
$ws.Range("H139").Value = 23450.572
$ws.Range("J139").Value = 23450.572
$ws.Range("L139").Value = 23450.572
$ws.Range("N139").Value = -33730.572
